# Update "想去人数" (number of people interested) values in column F
# for the "展览" and "全部类型" worksheets, reflecting newer scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    3  = 281
    5  = 840
    6  = 11
    7  = 296
    8  = 7927
    12 = 106
    15 = 19
    19 = 689
    20 = 21
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
